$d = $word.ActiveDocument
$d.TrackRevisions = $false

# ---------------------------------------------------------------------------
# Change 1: "GLOBAL SECURITY SYSTEMS" -> "GLOBAL SECURITY " + "SYSTEM" + " SOFTWARE"
#           (split into three runs with identical formatting)
# ---------------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("GLOBAL SECURITY SYSTEMS")
$r1 = $d.Range($find1.Start, $find1.End)
$full1 = "GLOBAL SECURITY SYSTEM SOFTWARE"
$r1.Text = $full1

$base1 = $r1.Start
$idxSystem = $full1.IndexOf("SYSTEM")
$sysStart = $base1 + $idxSystem
$sysEnd = $sysStart + "SYSTEM".Length
$sysRange = $d.Range($sysStart, $sysEnd)
# Force a run split without altering the final visible formatting: toggle Bold
# off then back on (text here is already bold), which leaves the run isolated
# from its neighbors without leaving any stray explicit overrides behind.
$sysRange.Bold = 0
$sysRange.Bold = 1

$softStart = $sysEnd
$softEnd = $base1 + $full1.Length
$softRange = $d.Range($softStart, $softEnd)
$softRange.Bold = 0
$softRange.Bold = 1

# ---------------------------------------------------------------------------
# Change 2: "..., shall be run only at United Nations Global Security Data
#           Centers around the world, and " ->
#           "..., shall be run only at United Nations " + "approved data
#           centers" + " around the world, and "
# ---------------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute([string][char]8221 + ", shall be run only at United Nations Global Security Data Centers around the world, and ")
$r2 = $d.Range($find2.Start, $find2.End)
$full2 = [string][char]8221 + ", shall be run only at United Nations approved data centers around the world, and "
$r2.Text = $full2

$base2 = $r2.Start
$needle2 = "approved data centers"
$idxApproved = $full2.IndexOf($needle2)
$apprStart = $base2 + $idxApproved
$apprEnd = $apprStart + $needle2.Length
$apprRange = $d.Range($apprStart, $apprEnd)
# This text is not bold in the source, so toggling Bold on then off isolates
# the run as its own element while leaving no explicit formatting residue.
$apprRange.Bold = 1
$apprRange.Bold = 0

# ---------------------------------------------------------------------------
# Change 3: merge the separate "," and " " runs (after "ANY CONFLICT") into a
#           single run containing ", "
# ---------------------------------------------------------------------------
$find3 = $d.Content
$find3.Find.Execute("ANY CONFLICT")
$commaPos = $find3.End
$spaceRange = $d.Range($commaPos + 1, $commaPos + 2)
$spaceRange.Delete()
$insertPoint = $d.Range($commaPos + 1, $commaPos + 1)
$insertPoint.InsertBefore(" ")

$d.Save()
